# "forgot to save WP2 update"
# Insert the missing WP2 task row ("Develop criteria for geometries/
# assemblies") that was dropped between the A/B date columns, shifting
# the later WP2/WP3/WP4/Report rows down by one, and correct the
# B-column (end) date on the row above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 18 (formatting is copied from row 17 above,
# which already carries the correct date-number-format styles for
# columns A and B).
$ws.Rows("18:18").Insert()

# Row 17 ("Create simplified CAD models of recommended components")
# actually finished on 14 Feb 2020, not 17 Feb 2020.
$ws.Range("B17").Value = 43875

# New row 18: the forgotten WP2 task. No start date (A18 left blank,
# matching the pattern used by other rows with only an end date),
# finished 17 Feb 2020, task name from the new shared string.
$ws.Range("B18").Value = 43878
$ws.Range("C18").Value = "Develop criteria for geometries/ assemblies"

# The sheet's data now spans one row further, so the autofilter and its
# backing _FilterDatabase defined name need to grow from C57/C58 to
# C58/C59 respectively (the filter range always included one blank
# trailing row beyond the data).
if ($ws.AutoFilterMode) {
  $ws.AutoFilterMode = $false
}
$ws.Range("A1:C58").AutoFilter()

foreach ($n in $wb.Names) {
  if ($n.Name -eq "Sheet1!_FilterDatabase") {
    $n.RefersTo = "=Sheet1!`$A`$1:`$C`$58"
  }
}

# Restore the cursor position that was active when the author saved.
$ws.Range("C9").Select()
